$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: SouthKorea / All - values updated
$ws.Cells.Item(2,4).Value = 0.02225547035720965
$ws.Cells.Item(2,5).Value = 0.1077507116422355
$ws.Cells.Item(2,6).Value = 0.0707498030982413
$ws.Cells.Item(2,7).Value = 0.03700090854399416
$ws.Cells.Item(2,8).Value = 0.6566063650062172
$ws.Cells.Item(2,9).Value = 0.343393634993783

# Row 3: China / All - values updated
$ws.Cells.Item(3,5).Value = 0.1071059312831127
$ws.Cells.Item(3,6).Value = 0.05290244791002885
$ws.Cells.Item(3,7).Value = 0.05420348337308388
$ws.Cells.Item(3,8).Value = 0.4939264079614041
$ws.Cells.Item(3,9).Value = 0.5060735920385959

# Row 4: Germany / All - values updated
$ws.Cells.Item(4,5).Value = 0.09193107182763709
$ws.Cells.Item(4,6).Value = 0.0486701250572603
$ws.Cells.Item(4,7).Value = 0.04326094677037679
$ws.Cells.Item(4,8).Value = 0.5294197499242979
$ws.Cells.Item(4,9).Value = 0.4705802500757021

# Row 5: USA / All - values updated
$ws.Cells.Item(5,4).Value = 0.0559832330827703
$ws.Cells.Item(5,5).Value = 0.07402294891667481
$ws.Cells.Item(5,6).Value = 0.06586423660656014
$ws.Cells.Item(5,7).Value = 0.008158712310114681
$ws.Cells.Item(5,8).Value = 0.8897813120185381
$ws.Cells.Item(5,9).Value = 0.110218687981462

# Row 6: previously Spain/All -> now USA/NYC with new values
$ws.Cells.Item(6,1).Value = "USA"
$ws.Cells.Item(6,2).Value = "NYC"
$ws.Cells.Item(6,4).Value = 0.07258978752642781
$ws.Cells.Item(6,5).Value = 0.0574163944730173
$ws.Cells.Item(6,6).Value = 0.06737503502022342
$ws.Cells.Item(6,7).Value = -0.009958640547206117
$ws.Cells.Item(6,8).Value = 0.8712250455686297
$ws.Cells.Item(6,9).Value = 0.1287749544313704

# Row 7: previously Italy/All -> now Spain/All with new values
$ws.Cells.Item(7,1).Value = "Spain"
$ws.Cells.Item(7,2).Value = "All"
$ws.Cells.Item(7,4).Value = 0.1040126934054379
$ws.Cells.Item(7,5).Value = 0.02599348859400723
$ws.Cells.Item(7,6).Value = 0.01257390976828207
$ws.Cells.Item(7,7).Value = 0.01341957882572516
$ws.Cells.Item(7,8).Value = 0.4837330596394426
$ws.Cells.Item(7,9).Value = 0.5162669403605574

# Row 8: new row Italy/All
$ws.Cells.Item(8,1).Value = "Italy"
$ws.Cells.Item(8,2).Value = "All"
$ws.Cells.Item(8,3).Value = 43943
$ws.Cells.Item(8,3).NumberFormat = $ws.Cells.Item(7,3).NumberFormat
$ws.Cells.Item(8,4).Value = 0.1300061819994451
$ws.Cells.Item(8,5).Value = 0
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = 0
